# Update the private health-insurance ("private Krankenkasse") sheet
# with the new contribution rates and the new entry date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Umlage U1 in Prozent: 1.6 -> 1.9
$ws.Range("B4").Value = 1.9

# Umlage U2 in Prozent: 0.44 -> 0.39
$ws.Range("B5").Value = 0.39

# Eintragungsdatum: 15.12.2023 -> 01.01.2024 (stored as text, like before)
$ws.Range("B7").Value = "01.01.2024"

# Leave the active cell on B8, like in the saved workbook
$ws.Range("B8").Select() | Out-Null
